$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.178113341331482
$ws.Range("B1").Value = 2.417391538619995
$ws.Range("D1").Value = 2.332413911819458
$ws.Range("E1").Value = 1.19875156879425
